# daily auto push: 2026-01-24 09:31 UTC
# A new 2026/01/24 17:00 observation (time=17, rank=138) was appended to the
# existing 2026/01/24 group. It belongs right after the last existing
# 2026/01/24 row (row 714), so insert a new row 715 and shift everything
# below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 715 (old "2026/12/29" row),
# pushing all rows from 715 downward one row further down.
$ws.Rows.Item(715).Insert()

# Fill in the new row with the inserted data point.
# Column A holds dates stored as plain text (not real Excel dates) in this
# sheet, so force a text entry via the leading apostrophe and then clear the
# style back to Normal so no extra text-format style gets attached to the
# cell (matching how the rest of the sheet's date cells are styled).
$ws.Cells.Item(715, 1).Value = "'2026/01/24"
$ws.Cells.Item(715, 1).Style = "Normal"
$ws.Cells.Item(715, 2).Value = "土"
$ws.Cells.Item(715, 3).Value = 17
$ws.Cells.Item(715, 4).Value = 138
